$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# Row 2 - "Preparacion de la Prueba" timing
$ws.Range("B2").Value = 0.003472222222222222
$ws.Range("D2").Value = 0.4791666666666667
$ws.Range("E2").Value = 0.4826388888888889

# Row 6 - "Crear interfaz Cola"
$ws.Range("A6").Value = "Crear interfaz Cola"
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 0.006944444444444444
$ws.Range("E6").Value = 0.4861111111111111
$ws.Range("F6").Value = 0.4895833333333333
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

# Row 7 - "Implementar Cola Estatica"
$ws.Range("A7").Value = "Implementar Cola Estatica"
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 41
$ws.Range("D7").Value = 0.013888888888888888
$ws.Range("E7").Value = 0.4909722222222222
$ws.Range("F7").Value = 0.5208333333333334
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0.020833333333333332

$ws.Range("B10").Select()

$wb.Save()
